# Apply the "Out of PO" roster update: rows 2-16 are reshuffled and
# "Luguentz Dort" (Oklahoma City Thunder) is replaced by "Grant Williams"
# (Charlotte Hornets).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Ochai Agbaji", "SG,SF", "Toronto Raptors"),
    @("Grant Williams", "PF,C", "Charlotte Hornets"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
